$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark (it will be re-added later
#    at its new location, right after the newly inserted "ICD10" text).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Expand "Includes " into the longer introductory sentence.
# ------------------------------------------------------------------
$find = $d.Content.Find
$find.Execute("Includes ", $true, $false, $false, $false, $false, $true, 1, $false,
              "Includes, for 1) new cases of dementia: ", 2) | Out-Null

# ------------------------------------------------------------------
# 3) Insert the new "Readcodes for primary care data, ICD10" text right
#    before the existing "Readcodes for conditions..." sentence, then
#    re-create the "_GoBack" bookmark immediately after "ICD10", and
#    finish with the new clause that leads back into the original text.
# ------------------------------------------------------------------
$find2 = $d.Content.Find
$find2.Execute("Readcodes for conditions", $true, $false, $false, $false, $false, $true, 1, $false,
               "", 0) | Out-Null
$insertPoint = $find2.Parent.Start

$rng = $d.Range($insertPoint, $insertPoint)
$rng.InsertBefore("Readcodes for primary care data, ICD10")

# Re-anchor on the text we just typed so we know exactly where "ICD10" ends.
$find3 = $d.Content.Find
$find3.Execute("ICD10", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$icd10End = $find3.Parent.End

$bmRange = $d.Range($icd10End, $icd10End)
$d.Bookmarks.Add("_GoBack", $bmRange)

$rng2 = $d.Range($icd10End, $icd10End)
$rng2.InsertAfter("- codes for Hospital Episode Statistics, product codes for medicinal products for dementia treatment; 2) for risk factors: ")

# ------------------------------------------------------------------
# 4) Split "Each list has an " so the sentence now starts with
#    "Where relevant, each list has an ".
# ------------------------------------------------------------------
$find4 = $d.Content.Find
$find4.Execute("Each list has an ", $true, $false, $false, $false, $false, $true, 1, $false,
               "Where relevant, each list has an ", 2) | Out-Null

Write-Host $d.Content.Text
